# Fixed the bubugs the cluster_prune_wrapperer, cluster_prune.m is deprecated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row above the old row 10 ("Dataset") for the new
#     "nn.cluster_base_quality_min = 0.3" sub-heading. This pushes all
#     the existing content (rows 10-25) down by one row (11-26), and
#     Excel automatically fixes up the shared-formula ranges.
$ws.Rows("10:10").Insert()
$ws.Range("A10").Value2 = "nn.cluster_base_quality_min = 0.3"

# --- Append the new "base quality 0.2" experiment block starting at
#     row 28 (row 27 left blank as a separator, matching existing style).

# Fill in the data tables first...
$ws.Range("A29").Value2 = "nn.cluster_base_quality_min = 0.2"

# SVHN table
$ws.Range("A30").Value2 = "SVHN (1024-1200-1200-10)"
$ws.Range("B30").Value2 = "Unclustered"
$ws.Range("C30").Value2 = "Clustered"
$ws.Range("D30").Value2 = "Ratio"

$ws.Range("A31").Value2 = "Layer 1"
$ws.Range("B31").Value2 = 15174
$ws.Range("C31").Value2 = 8712
$ws.Range("D31").Formula = "=B31/C31"

$ws.Range("A32").Value2 = "Layer 2"
$ws.Range("B32").Value2 = 13733
$ws.Range("C32").Value2 = 58605
$ws.Range("D32").Formula = "=B32/C32"

$ws.Range("A33").Value2 = "Layer 3"
$ws.Range("B33").Value2 = 542
$ws.Range("C33").Value2 = 785
$ws.Range("D33").Formula = "=B33/C33"

# MNIST table
$ws.Range("A35").Value2 = "MNIST (784-1200-1200-10)"
$ws.Range("B35").Value2 = "Unclustered"
$ws.Range("C35").Value2 = "Clustered"

$ws.Range("A36").Value2 = "Layer 1"
$ws.Range("B36").Value2 = 13414
$ws.Range("C36").Value2 = 11746
$ws.Range("D36").Formula = "=B36/C36"

$ws.Range("A37").Value2 = "Layer 2"
$ws.Range("B37").Value2 = 23899
$ws.Range("C37").Value2 = 8977
$ws.Range("D37").Formula = "=B37/C37"

$ws.Range("A38").Value2 = "Layer 3"
$ws.Range("B38").Value2 = 496
$ws.Range("C38").Value2 = 387
$ws.Range("D38").Formula = "=B38/C38"

# CIFAR-10 table
$ws.Range("A40").Value2 = "CIFAR-10 (1024-1200-1200-10)"
$ws.Range("B40").Value2 = "Unclustered"
$ws.Range("C40").Value2 = "Clustered"

$ws.Range("A41").Value2 = "Layer 1"
$ws.Range("B41").Value2 = 31179
$ws.Range("C41").Value2 = 48248
$ws.Range("D41").Formula = "=B41/C41"

$ws.Range("A42").Value2 = "Layer 2"
$ws.Range("B42").Value2 = 25075
$ws.Range("C42").Value2 = 45235
$ws.Range("D42").Formula = "=B42/C42"

$ws.Range("A43").Value2 = "Layer 3"
$ws.Range("B43").Value2 = 502
$ws.Range("C43").Value2 = 849
$ws.Range("D43").Formula = "=B43/C43"

# Conclusions for this experiment
$ws.Range("A44").Value2 = "Conslusion: Lower base quality for xbar utilization improved the ratio of unclustered to clustered"
$ws.Range("A45").Value2 = "However, #unclustered crossbars are still high."
$ws.Range("A46").Value2 = "Lower accuracy in SVHN for transformers than just pruning."
$ws.Range("A47").Value2 = "Lower pruning in SVHN and CIFAR-10"

# ... and type the heading for this block last.
$ws.Range("A28").Value2 = "Trying a lower base quality to reduce the number of unclustered synapses and resulting crossbars"

# --- Apply formatting to match the rest of the sheet ---
# A1-style (bold) headings
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A28").PasteSpecial(-4122) | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null

# A2-style (plain, left-aligned) table rows/headers
$ws.Range("A11:D11").Copy() | Out-Null
$ws.Range("A30:D30").PasteSpecial(-4122) | Out-Null

$ws.Range("A13:D15").Copy() | Out-Null
$ws.Range("A31:D33").PasteSpecial(-4122) | Out-Null

$ws.Range("A17:C17").Copy() | Out-Null
$ws.Range("A35:C35").PasteSpecial(-4122) | Out-Null

$ws.Range("A18:D20").Copy() | Out-Null
$ws.Range("A36:D38").PasteSpecial(-4122) | Out-Null

$ws.Range("A22:C22").Copy() | Out-Null
$ws.Range("A40:C40").PasteSpecial(-4122) | Out-Null

$ws.Range("A23:D25").Copy() | Out-Null
$ws.Range("A41:D43").PasteSpecial(-4122) | Out-Null

# A3-style (bold red) conclusion line
$ws.Range("A26").Copy() | Out-Null
$ws.Range("A44").PasteSpecial(-4122) | Out-Null

# remaining plain text lines (A2-style)
$ws.Range("A20").Copy() | Out-Null
$ws.Range("A45").PasteSpecial(-4122) | Out-Null
$ws.Range("A46").PasteSpecial(-4122) | Out-Null
$ws.Range("A47").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Append the second "iso-accuracy pruning" block, which duplicates
#     the SVHN and CIFAR-10 tables from the block above ---

# SVHN table (same data as the block above)
$ws.Range("A50").Value2 = "SVHN (1024-1200-1200-10)"
$ws.Range("B50").Value2 = "Unclustered"
$ws.Range("C50").Value2 = "Clustered"
$ws.Range("D50").Value2 = "Ratio"

$ws.Range("A51").Value2 = "Layer 1"
$ws.Range("B51").Value2 = 15174
$ws.Range("C51").Value2 = 8712
$ws.Range("D51").Formula = "=B51/C51"

$ws.Range("A52").Value2 = "Layer 2"
$ws.Range("B52").Value2 = 13733
$ws.Range("C52").Value2 = 58605
$ws.Range("D52").Formula = "=B52/C52"

$ws.Range("A53").Value2 = "Layer 3"
$ws.Range("B53").Value2 = 542
$ws.Range("C53").Value2 = 785
$ws.Range("D53").Formula = "=B53/C53"

# CIFAR-10 table (same data as the block above)
$ws.Range("A55").Value2 = "CIFAR-10 (1024-1200-1200-10)"
$ws.Range("B55").Value2 = "Unclustered"
$ws.Range("C55").Value2 = "Clustered"

$ws.Range("A56").Value2 = "Layer 1"
$ws.Range("B56").Value2 = 31179
$ws.Range("C56").Value2 = 48248
$ws.Range("D56").Formula = "=B56/C56"

$ws.Range("A57").Value2 = "Layer 2"
$ws.Range("B57").Value2 = 25075
$ws.Range("C57").Value2 = 45235
$ws.Range("D57").Formula = "=B57/C57"

$ws.Range("A58").Value2 = "Layer 3"
$ws.Range("B58").Value2 = 502
$ws.Range("C58").Value2 = 849
$ws.Range("D58").Formula = "=B58/C58"

# Formatting for this block (copied from the corresponding block above).
$ws.Range("A30:D33").Copy() | Out-Null
$ws.Range("A50:D53").PasteSpecial(-4122) | Out-Null

$ws.Range("A40:C40").Copy() | Out-Null
$ws.Range("A55:C55").PasteSpecial(-4122) | Out-Null

$ws.Range("A41:D43").Copy() | Out-Null
$ws.Range("A56:D58").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Heading for this block, typed last (as with the block above).
$ws.Range("A49").Value2 = "Trying a iso-accuracy pruning for prunemodes - 1 & 2, by forced prune slowdown in mode 1 "
$ws.Range("A28").Copy() | Out-Null
$ws.Range("A49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Put the selection roughly where the author left it.
$ws.Range("B23").Select()
